# update made january 24th
# Append three new churn records (rows 13-15) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking CVR number as TEXT (not auto-converted to a
# number) by staging it through a text formula in a scratch cell far outside
# the used range, copying the resulting (string-typed) cell into place, then
# clearing the scratch cell again. A direct `.Value = "34391513"` assignment
# would be auto-coerced to a Number by Excel, which is not what we want here
# (column A stores CVR numbers as text, matching the rest of the column).
function Set-TextValue {
    param($range, [string]$text)

    $scratch = $ws.Range("ZZ1")
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy($range) | Out-Null
    $scratch.ClearContents() | Out-Null
}

# --- Row 13 ---------------------------------------------------------------
Set-TextValue $ws.Range("A13") "34391513"
$ws.Range("B13").Value = 2023
$ws.Range("C13").Value = 102035.85
$ws.Range("D13").Value = "Visma Løn"
$ws.Range("E13").Value = 45196
$ws.Range("E13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G13").Value = "DataLøn"
$ws.Range("H13").Value = "2023Q3"
$ws.Range("I13").Value = "100000-120000"

# --- Row 14 ---------------------------------------------------------------
Set-TextValue $ws.Range("A14") "10613779"
$ws.Range("B14").Value = 2023
$ws.Range("C14").Value = 116715
$ws.Range("D14").Value = "Visma Løn"
$ws.Range("E14").Value = 45212
$ws.Range("E14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H14").Value = "2023Q4"
$ws.Range("I14").Value = "100000-120000"

# --- Row 15 ---------------------------------------------------------------
Set-TextValue $ws.Range("A15") "66328511"
$ws.Range("B15").Value = 2023
$ws.Range("C15").Value = 109559.83
$ws.Range("D15").Value = "Visma Løn og HR"
$ws.Range("E15").Value = 45245
$ws.Range("E15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H15").Value = "2023Q4"
$ws.Range("I15").Value = "100000-120000"
